{"js": "// \"Small update to history\"\n// Two small textual corrections in the \"One day, at lunch ...\" paragraph:\n//   1. \"One day, at lunch\"  ->  \"One day in 1993, at lunch\"\n//   2. \"He knew nothing about playing cards, let alone bridge\"\n//         -> \"He knew little about card games, and nothing about bridge\"\n\nconst body = context.document.body;\n\n// 1) Insert \" in 1993\" right before \", at lunch\".\nconst lunchResults = body.search(\", at lunch\", { matchCase: true });\nlunchResults.load(\"text\");\nawait context.sync();\n\nif (lunchResults.items.length > 0) {\n  const lunchRange = lunchResults.items[0];\n  // Put the new text immediately before the matched range.\n  lunchRange.insertText(\" in 1993\", \"Before\");\n}\n\n// 2) Swap out the \"nothing about playing cards, let alone\" phrasing.\nconst knewResults = body.search(\n  \"He knew nothing about playing cards, let alone \",\n  { matchCase: true }\n);\nknewResults.load(\"text\");\nawait context.sync();\n\nif (knewResults.items.length > 0) {\n  knewResults.items[0].insertText(\n    \"He knew little about card games, and nothing about \",\n    \"Replace\"\n  );\n}\n\nawait context.sync();\n", "ps1": "# \"Small update to history\"\n# Two small textual corrections in the \"One day, at lunch ...\" paragraph:\n#   1. \"One day, at lunch\"  ->  \"One day in 1993, at lunch\"\n#   2. \"He knew nothing about playing cards, let alone bridge\"\n#         -> \"He knew little about card games, and nothing about bridge\"\n\n$d = $word.ActiveDocument\n\n# 1) Insert \" in 1993\" right before \", at lunch\".\n$rng1 = $d.Content\n$rng1.Find.MatchCase = $true\n$rng1.Find.Forward = $true\n$found1 = $rng1.Find.Execute(\", at lunch\")\nif ($found1) {\n    $rng1.Collapse(1)   # wdCollapseStart\n    $rng1.InsertBefore(\" in 1993\")\n}\n\n# 2) Swap out the \"nothing about playing cards, let alone\" phrasing.\n$rng2 = $d.Content\n$rng2.Find.MatchCase = $true\n$rng2.Find.Forward = $true\n$found2 = $rng2.Find.Execute(\"He knew nothing about playing cards, let alone \")\nif ($found2) {\n    $rng2.Text = \"He knew little about card games, and nothing about \"\n}\n"}
